# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" sheets, which carry duplicate copies of the data.
$wb = $excel.ActiveWorkbook

# row -> new F value
$updates = @{
    4  = 647
    6  = 23
    7  = 11706
    13 = 825
    14 = 13406
    15 = 13254
    17 = 148
    20 = 261
    22 = 47
    23 = 143
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
